$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2, column E (Authors) - updated author list (Fritz Francois flag changed 1 -> 0, extra spacing)
$ws.Range("E2").Value = "[Jennifer%Lighter%Jennifer.Lighter@nyumc.org%1,        Michael%Phillips%NULL%1,        Sarah%Hochman%NULL%1,        Stephanie%Sterling%NULL%1,        Diane%Johnson%NULL%1,        Fritz%Francois%NULL%0,        Anna%Stachel%NULL%1]"

# Row 3 - new CrossRef-sourced record replacing the previously "not found" placeholder data
$ws.Range("C3").Value = "`"Early antiviral treatment contributes to alleviate the severity and improve the prognosis of patients with novel coronavirus disease (COVID\u201019)`""
$ws.Range("E3").Value = "[J.%Wu%xref no email%1,  W.%Li%xref no email%1,  X.%Shi%xref no email%1,  Z.%Chen%xref no email%2,  B.%Jiang%xref no email%1,  J.%Liu%xref no email%0,  D.%Wang%xref no email%1,  C.%Liu%xref no email%1,  Y.%Meng%xref no email%1,  L.%Cui%xref no email%1,  J.%Yu%xref no email%1,  H.%Cao%xref no email%1,  L.%Li%xref no email%1]"
$ws.Range("F3").Value = "10.1111/joim.13063"
$ws.Range("G3").Value = "CROSSREF"

# H3 ("Date Accepted") must stay plain text "2023-05-09", not be auto-converted to
# a date serial number by Excel's smart entry parsing. Route it through a text
# formula + paste-values so it lands as a literal shared string with no style churn.
$ws.Range("H3").Formula = "=T(""2023-05-09"")"
$ws.Range("H3").Copy() | Out-Null
$ws.Range("H3").PasteSpecial(-4163) | Out-Null
